# "Revision antes de 3ra produccion SETI" / "Se actualiza la columna NEW"
#
# 1) BATT_CONNECTOR: the "NEW" column (B7:B9) flips from YES to NO ahead of
#    the 3rd SETI production batch.
# 2) _HISTORY: a new version-4 row is logged documenting this change.

$wb = $excel.ActiveWorkbook

$wsBC = $wb.Worksheets.Item("BATT_CONNECTOR")
$wsBC.Range("B7").Value = "NO"
$wsBC.Range("B8").Value = "NO"
$wsBC.Range("B9").Value = "NO"

$wsH = $wb.Worksheets.Item("_HISTORY")
$wsH.Range("A7").Value = 4
$newDate = Get-Date -Year 2023 -Month 9 -Day 26
# Set the number format first so Excel doesn't provision a transient
# "m/d/yyyy" custom format as a side effect of assigning a date Value.
$wsH.Range("B7").NumberFormat = "d-mmm-yy"
$wsH.Range("B7").Value = $newDate.Date
$wsH.Range("C7").Value = "DGB"
$wsH.Range("D7").Value = "Se actualiza columna NEW para 3ra produccion de SETI"

# Leave the selection where the author left it when saving.
$wsH.Activate() | Out-Null
$wsH.Range("B8").Select() | Out-Null

$wsBC.Activate() | Out-Null
$wsBC.Range("B10").Select() | Out-Null
